$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.75
$ws.Range("AC3").Value = 21
$ws.Range("AN3").Value = 23

# Row 4
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.93
$ws.Range("AS4").Value = 1.47

# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 2.8
$ws.Range("L5").Value = 3.25
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 15
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 5
$ws.Range("Q5").Value = 1.6
$ws.Range("R5").Value = 2.3
$ws.Range("S5").Value = 2.38
$ws.Range("T5").Value = 1.53
$ws.Range("U5").Value = 1.3
$ws.Range("V5").Value = 3.4
$ws.Range("W5").Value = 1.5
$ws.Range("X5").Value = 2.5
$ws.Range("Y5").Value = 11
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 15
$ws.Range("AF5").Value = 7
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 34
$ws.Range("AI5").Value = 101
$ws.Range("AJ5").Value = 13
$ws.Range("AO5").Value = 23
$ws.Range("AR5").Value = 1.98
$ws.Range("AS5").Value = 1.88

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 2.87
$ws.Range("I8").Value = 3.7
$ws.Range("J8").Value = 2.82
$ws.Range("K8").Value = 1.93
$ws.Range("L8").Value = 4.15
$ws.Range("N8").Value = 5.6
$ws.Range("O8").Value = 1.47
$ws.Range("P8").Value = 2.52
$ws.Range("Q8").Value = 2.35
$ws.Range("S8").Value = 4.15
$ws.Range("T8").Value = 1.19
$ws.Range("Y8").Value = 5.8
$ws.Range("Z8").Value = 9.25
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 21
$ws.Range("AD8").Value = 37
$ws.Range("AE8").Value = 5.6
$ws.Range("AF8").Value = 5.6
$ws.Range("AG8").Value = 15.5
$ws.Range("AJ8").Value = 9
$ws.Range("AL8").Value = 12.5
$ws.Range("AN8").Value = 37
$ws.Range("AO8").Value = 45

# Row 9
$ws.Range("J9").Value = 3.75
$ws.Range("U9").Value = 1.36
$ws.Range("V9").Value = 3
$ws.Range("W9").Value = 1.67
$ws.Range("X9").Value = 2.1
$ws.Range("AE9").Value = 12
$ws.Range("AK9").Value = 11
$ws.Range("AM9").Value = 19

# Row 15
$ws.Range("I15").Value = 1.85
$ws.Range("L15").Value = 2.63
$ws.Range("Y15").Value = 9.5
$ws.Range("AJ15").Value = 6
$ws.Range("AK15").Value = 8

# Row 16
$ws.Range("G16").Value = 1.95
$ws.Range("I16").Value = 3.8
$ws.Range("J16").Value = 2.75
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 4.5
$ws.Range("O16").Value = 1.4
$ws.Range("P16").Value = 2.75
$ws.Range("Y16").Value = 6.5
$ws.Range("AC16").Value = 19
$ws.Range("AI16").Value = 800
$ws.Range("AK16").Value = 19
$ws.Range("AN16").Value = 34

# Row 17
$ws.Range("H17").Value = 3.1
$ws.Range("K17").Value = 2
$ws.Range("N17").Value = 7.5
$ws.Range("S17").Value = 4.33
$ws.Range("T17").Value = 1.2
$ws.Range("AD17").Value = 34
$ws.Range("AE17").Value = 7.5
$ws.Range("AI17").Value = 900
$ws.Range("AP17").Value = 1.75
$ws.Range("AQ17").Value = 2.05

# Row 18
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13
$ws.Range("O18").Value = 1.29
$ws.Range("P18").Value = 3.5
$ws.Range("Q18").Value = 1.9
$ws.Range("R18").Value = 1.95

# Row 21
$ws.Range("G21").Value = 1.75
$ws.Range("H21").Value = 4.1
$ws.Range("I21").Value = 3.85
$ws.Range("J21").Value = 2.18
$ws.Range("K21").Value = 2.52
$ws.Range("L21").Value = 3.95
$ws.Range("Q21").Value = 1.38
$ws.Range("R21").Value = 2.8
$ws.Range("S21").Value = 1.93
$ws.Range("T21").Value = 1.78
$ws.Range("U21").Value = 1.23
$ws.Range("V21").Value = 3.75
$ws.Range("X21").Value = 2.7
$ws.Range("Z21").Value = 12.5
$ws.Range("AB21").Value = 17
$ws.Range("AC21").Value = 11.75
$ws.Range("AD21").Value = 16
$ws.Range("AF21").Value = 9
$ws.Range("AG21").Value = 11.75
$ws.Range("AH21").Value = 32
$ws.Range("AJ21").Value = 19.5
$ws.Range("AK21").Value = 28
$ws.Range("AL21").Value = 13.5
$ws.Range("AM21").Value = 60
$ws.Range("AN21").Value = 28
$ws.Range("AO21").Value = 25

# Row 22
$ws.Range("G22").Value = 1.75
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 5.1
$ws.Range("J22").Value = 2.4
$ws.Range("K22").Value = 1.95
$ws.Range("L22").Value = 5.5
$ws.Range("M22").Value = 1.11
$ws.Range("N22").Value = 5.6
$ws.Range("O22").Value = 1.47
$ws.Range("P22").Value = 2.5
$ws.Range("Q22").Value = 2.37
$ws.Range("R22").Value = 1.52
$ws.Range("S22").Value = 4.25
$ws.Range("T22").Value = 1.18
$ws.Range("U22").Value = 1.52
$ws.Range("V22").Value = 2.37
$ws.Range("W22").Value = 2.15
$ws.Range("X22").Value = 1.62
$ws.Range("Z22").Value = 7
$ws.Range("AA22").Value = 8.75
$ws.Range("AB22").Value = 14
$ws.Range("AC22").Value = 17.5
$ws.Range("AD22").Value = 40
$ws.Range("AE22").Value = 5.6
$ws.Range("AF22").Value = 6.2
$ws.Range("AG22").Value = 19
$ws.Range("AH22").Value = 120
$ws.Range("AJ22").Value = 11
$ws.Range("AK22").Value = 30
$ws.Range("AL22").Value = 17
$ws.Range("AM22").Value = 110
$ws.Range("AN22").Value = 65
$ws.Range("AO22").Value = 70

# Row 24
$ws.Range("G24").Value = 2.3
$ws.Range("H24").Value = 3.25
$ws.Range("J24").Value = 2.9
$ws.Range("K24").Value = 2.12
$ws.Range("M24").Value = 1.05
$ws.Range("N24").Value = 7.6
$ws.Range("O24").Value = 1.27
$ws.Range("P24").Value = 3.45
$ws.Range("Q24").Value = 1.8
$ws.Range("R24").Value = 1.91
$ws.Range("S24").Value = 2.9
$ws.Range("T24").Value = 1.36
$ws.Range("V24").Value = 2.82
$ws.Range("X24").Value = 2.12
$ws.Range("AD24").Value = 25
$ws.Range("AE24").Value = 7.6
$ws.Range("AG24").Value = 12.5
$ws.Range("AH24").Value = 50
$ws.Range("AI24").Value = 350
$ws.Range("AJ24").Value = 9.75
$ws.Range("AK24").Value = 15.5
$ws.Range("AO24").Value = 29

# Row 28
$ws.Range("H28").Value = 3.7
$ws.Range("I28").Value = 2.1
$ws.Range("J28").Value = 3.4
$ws.Range("K28").Value = 2.32
$ws.Range("L28").Value = 2.62
$ws.Range("M28").Value = 1.03
$ws.Range("N28").Value = 9
$ws.Range("P28").Value = 4.2
$ws.Range("U28").Value = 1.3
$ws.Range("V28").Value = 3.2
$ws.Range("W28").Value = 1.52
$ws.Range("X28").Value = 2.37
$ws.Range("Y28").Value = 13
$ws.Range("Z28").Value = 18.5
$ws.Range("AC28").Value = 22
$ws.Range("AE28").Value = 9
$ws.Range("AF28").Value = 7.5
$ws.Range("AG28").Value = 11.75
$ws.Range("AH28").Value = 40
$ws.Range("AI28").Value = 250
$ws.Range("AJ28").Value = 10.5
$ws.Range("AK28").Value = 12.5
$ws.Range("AM28").Value = 21
$ws.Range("AO28").Value = 20

# Row 29
$ws.Range("G29").Value = 2.65
$ws.Range("I29").Value = 2.4
$ws.Range("J29").Value = 3.1
$ws.Range("K29").Value = 2.25
$ws.Range("L29").Value = 2.92
$ws.Range("O29").Value = 1.18
$ws.Range("P29").Value = 4.3
$ws.Range("Q29").Value = 1.55
$ws.Range("R29").Value = 2.3
$ws.Range("S29").Value = 2.3
$ws.Range("T29").Value = 1.55
$ws.Range("U29").Value = 1.3
$ws.Range("V29").Value = 3.2
$ws.Range("W29").Value = 1.47
$ws.Range("X29").Value = 2.52
$ws.Range("Y29").Value = 13
$ws.Range("Z29").Value = 17
$ws.Range("AA29").Value = 9.75
$ws.Range("AB29").Value = 32
$ws.Range("AC29").Value = 19
$ws.Range("AD29").Value = 21
$ws.Range("AF29").Value = 7.2
$ws.Range("AH29").Value = 35
$ws.Range("AJ29").Value = 11.5
$ws.Range("AK29").Value = 14.5
$ws.Range("AL29").Value = 9.25
$ws.Range("AM29").Value = 27
$ws.Range("AN29").Value = 17
